$d = $word.ActiveDocument

$replacements = @(
    @("2023-09-23 Saturday", "2023-09-24 Sunday"),
    @("63+12=", "83-82="),
    @("86-24=", "71+18="),
    @("47-27=", "81-11="),
    @("82-49=", "35-5="),
    @("3+48=", "27+34="),
    @("73-47=", "50-27="),
    @("2+7=", "76+4="),
    @("54+10=", "36+41="),
    @("85-27=", "2+25="),
    @("80-76=", "57+34="),
    @("17-13=", "33+62="),
    @("49+35=", "81-57="),
    @("68-63=", "14+49="),
    @("62-9=", "55-0="),
    @("61-55=", "28+67="),
    @("29+0=", "10+26="),
    @("68-35=", "67-49="),
    @("15-6=", "32+58="),
    @("75-57=", "70-27="),
    @("66-7=", "51-48="),
    @("13+22=", "16+17="),
    @("28+71=", "29+58="),
    @("11+54=", "95-70="),
    @("95+0=", "29+20="),
    @("57-12=", "98-40="),
    @("24-15=", "49+12="),
    @("62+26=", "15+23="),
    @("35+22=", "2+63="),
    @("79-76=", "39-33="),
    @("45+8=", "17-10="),
    @("94-6=", "49+11="),
    @("89-74=", "90-3="),
    @("13+33=", "16-5="),
    @("9+8=", "51-48="),
    @("25+73=", "93-36="),
    @("92-57=", "97-60="),
    @("37+25=", "83-45="),
    @("10+64=", "49+41="),
    @("70-64=", "93-7="),
    @("84-23=", "93-79="),
    @("67-27=", "82-62="),
    @("31+34=", "69-40="),
    @("97-4=", "14+46="),
    @("2-0=", "57-52="),
    @("98-70=", "2+18="),
    @("52+18=", "25+43="),
    @("60+20=", "73+2="),
    @("81-18=", "79-70="),
    @("1+35=", "24+4="),
    @("99-81=", "53-9="),
    @("67-57=", "96-52="),
    @("41+36=", "44-36="),
    @("51+39=", "87-38="),
    @("67-23=", "27+21="),
    @("9+78=", "39+1="),
    @("76+20=", "11-7="),
    @("24+5=", "86-81="),
    @("43+21=", "72+23="),
    @("47+12=", "34+59="),
    @("10+38=", "8+41="),
    @("33+2=", "34-8="),
    @("91-3=", "87-62="),
    @("26+41=", "39+46="),
    @("24+44=", "63+9="),
    @("57+9=", "59-54="),
    @("20+20=", "70-59="),
    @("32+54=", "78+16="),
    @("80-79=", "40-32="),
    @("40-1=", "84-51="),
    @("13+10=", "10+19="),
    @("44+10=", "93-32="),
    @("28+31=", "63-2="),
    @("83+1=", "58+18="),
    @("16+36=", "78-14="),
    @("17+9=", "2+62="),
    @("23+0=", "39-8="),
    @("80-4=", "21+37="),
    @("12-4=", "96-92="),
    @("63-11=", "46+43="),
    @("19+21=", "36+2="),
    @("93-34=", "57-36="),
    @("48+30=", "63-10="),
    @("12+86=", "93-23="),
    @("59-34=", "67+2="),
    @("4+57=", "66+9="),
    @("66-52=", "43-33="),
    @("25+14=", "64-27="),
    @("7+63=", "10-8="),
    @("74-56=", "45-45="),
    @("99-21=", "92-76="),
    @("27-19=", "5+53="),
    @("45+46=", "91-89="),
    @("85-9=", "79-54="),
    @("60-1=", "44-0="),
    @("34+6=", "4+66="),
    @("91-42=", "98-34="),
    @("76-43=", "15+7="),
    @("17+69=", "25+28="),
    @("54+0=", "17-2="),
    @("84-4=", "22+41="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
